# Auto-generated: update cryptos list with latest price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting is preserved for price-like strings (avoid numeric/date auto-conversion)
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "65.450.75"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "3.265.32"
$ws.Range("E3").Value = "  -1.63%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "575.48"
$ws.Range("E5").Value = "  +2.84%  "

$ws.Range("D6").Value = "181.08"
$ws.Range("E6").Value = "  -3.59%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "3.262.90"
$ws.Range("E8").Value = "  -1.47%  "

$ws.Range("D9").Value = "0.567"
$ws.Range("E9").Value = "  -3.18%  "

$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  -6.99%  "

$ws.Range("D11").Value = "0.565"
$ws.Range("E11").Value = "  -3.93%  "

$ws.Range("D12").Value = "45.80"
$ws.Range("E12").Value = "  -4.06%  "

$ws.Range("D13").Value = "0.0000261"
$ws.Range("E13").Value = "  -4.20%  "

$ws.Range("D14").Value = "3.793.64"
$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").Value = "8.31"
$ws.Range("E15").Value = "  -4.14%  "

$ws.Range("D16").Value = "611.27"
$ws.Range("E16").Value = "  -3.51%  "

$ws.Range("D17").Value = "65.593.03"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.276.84"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "17.52"
$ws.Range("E20").Value = "  -3.48%  "

$ws.Range("D21").Value = "10.78"
$ws.Range("E21").Value = "  -2.07%  "

$ws.Range("D22").Value = "0.879"
$ws.Range("E22").Value = "  -3.44%  "

$ws.Range("D23").Value = "18.12"
$ws.Range("E23").Value = "  -0.93%  "

$ws.Range("D24").Value = "4.90"
$ws.Range("E24").Value = "  -0.95%  "

$ws.Range("D25").Value = "97.86"
$ws.Range("E25").Value = "  -5.07%  "

$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -1.39%  "

$ws.Range("D28").Value = "9.31"
$ws.Range("E28").Value = "  -3.35%  "

$ws.Range("D29").Value = "30.34"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").Value = "8.29"
$ws.Range("E30").Value = "  -4.74%  "

$ws.Range("D31").Value = "6.39"
$ws.Range("E31").Value = "  -0.69%  "

$ws.Range("D32").Value = "3.70"
$ws.Range("E32").Value = "  -6.84%  "

$ws.Range("D33").Value = "543.20"
$ws.Range("E33").Value = "  -1.89%  "

$ws.Range("D34").Value = "10.74"
$ws.Range("E34").Value = "  -3.28%  "

$ws.Range("D35").Value = "3.751.17"
$ws.Range("E35").Value = "  -2.17%  "

$ws.Range("E36").Value = "  -3.20%  "

$ws.Range("D37").Value = "0.997"
$ws.Range("E37").Value = "  -0.33%  "

$ws.Range("D38").Value = "55.60"
$ws.Range("E38").Value = "  -3.67%  "

$ws.Range("D39").Value = "0.126"
$ws.Range("E39").Value = "  -1.61%  "

$ws.Range("D40").Value = "32.19"
$ws.Range("E40").Value = "  -5.43%  "

$ws.Range("D41").Value = "3.37"
$ws.Range("E41").Value = "  +3.86%  "

$ws.Range("D42").Value = "3.11"
$ws.Range("E42").Value = "  -5.91%  "

$ws.Range("D43").Value = "0.0₃0669"
$ws.Range("E43").Value = "  -9.09%  "

$ws.Range("D44").Value = "2.55"
$ws.Range("E44").Value = "  -6.05%  "

$ws.Range("D45").Value = "0.326"
$ws.Range("E45").Value = "  -3.14%  "

$ws.Range("D46").Value = "0.0402"
$ws.Range("E46").Value = "  -4.66%  "

$ws.Range("D47").Value = "2.97"
$ws.Range("E47").Value = "  -7.98%  "

$ws.Range("E48").Value = "  +0.31%  "

$ws.Range("D49").Value = "0.125"
$ws.Range("E49").Value = "  -3.21%  "

$ws.Range("D50").Value = "2.48"
$ws.Range("E50").Value = "  -5.15%  "

$ws.Range("D51").Value = "127.80"
$ws.Range("E51").Value = "  +4.45%  "

